# 15.07.2025 - maç sonuçları eklendi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Yarı Final 1 (row 15): Ajans Of vs Fortuna United -> 2-2, penaltılar 11-10
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 2
$ws.Range("I15").Value = 11
$ws.Range("J15").Value = 10

# Yarı Final 2 (row 18): Araklı 1961 Spor vs Hubuş FK -> 2-2, penaltılar 12-11
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 2
$ws.Range("I18").Value = 12
$ws.Range("J18").Value = 11

# Final (row 21): Ajans Of vs Araklı 1961 Spor
$ws.Range("A21").Value = "Ajans Of"
$ws.Range("B21").Value = "Araklı 1961 Spor"

# 3. lük Maçı (row 24): Hubuş FK vs Fortuna United
$ws.Range("A24").Value = "Hubuş FK"
$ws.Range("B24").Value = "Fortuna United"

# Update selection to reflect last-used cell
$ws.Range("D25").Select()
